$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LLBV3 "Header / Function" notes in columns D/E were reassigned to
# different physical pins:
#   - BUZZER / "on-board buzzer" moves from row 53 (PG1/RD, Digital pin 40)
#     to row 51 (PD7/T0, Digital pin 38).
#   - BRAKE_ON / "on-board relay" moves from row 52 (PG0/WR, Digital pin 41)
#     to row 71 (PG2/ALE, Digital pin 39); the relay note is retyped with a
#     space instead of a hyphen ("on board relay").
$ws.Range("D51").Value = "BUZZER"
$ws.Range("E51").Value = "on-board buzzer"

$ws.Range("D52").ClearContents()
$ws.Range("E52").ClearContents()

$ws.Range("D53").ClearContents()
$ws.Range("E53").ClearContents()

$ws.Range("D71").Value = "BRAKE_ON"
$ws.Range("E71").Value = "on board relay"

# Scroll the view to where the edits were made and leave the selection on
# the cell below the newly-entered relay note.
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 2
$ws.Range("D72").Select()
